# persAlignments.xlsx update:
#  - new 4th glm model (age + age^2 + individual items) written to columns S (new
#    dense item model) and the existing "A" model re-emitted in column X with an
#    age^2 row inserted.
#  - bold "CORRECT FOR MEDICATION" header label added at N1:O1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Bold header label in N1:O1
# ---------------------------------------------------------------------------
$ws.Range("N1").Value = "CORRECT FOR MEDICATION"
$ws.Range("N1:O1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) New column S: full item-level glm() printout (intercept, sex, age, age2,
#    then every personality item) in the same fixed-width "Lucida Console"
#    style used throughout the sheet.
# ---------------------------------------------------------------------------
$sVals = @{
    1 = "(Intercept)    0.0628444918"
    2 = "sex            .           "
    3 = "age            0.0430165272"
    4 = "age2           0.1117059672"
    5 = "Depressed      .           "
    6 = "Fearful        .           "
    7 = "Persistent     .           "
    8 = "Cautious       .           "
    9 = "Stable         .           "
    10 = "Autistic      -0.0473130008"
    11 = "Stingy         .           "
    12 = "Jealous        .           "
    13 = "Reckless       .           "
    14 = "Sociable       .           "
    15 = "Timid         -0.0142789320"
    16 = "Sympathetic    .           "
    17 = "Playful        .           "
    18 = "Solitary       .           "
    19 = "Active        -0.0974213317"
    20 = "Helpful        .           "
    21 = "Bullying       .           "
    22 = "Aggressive     .           "
    23 = "Manipiulative  .           "
    24 = "Gentle         .           "
    25 = "Affectionate   0.0266638768"
    26 = "Excitable     -0.0000547059"
    27 = "Impulsive      .           "
    28 = "Inquisitve     .           "
    29 = "Submissive     .           "
    30 = "Dependent      .           "
    31 = "Irritible      .           "
    32 = "Predictable    .           "
    33 = "Decisive       0.0499958689"
    34 = "Independent    .           "
    35 = "Sensitive      .           "
    36 = "Defiant        .           "
    37 = "Intelligent    0.0268132986"
    38 = "Protective     .           "
    39 = "Inventive      .           "
    40 = "Clumsy         .           "
    41 = "Erratic       -0.0282228804"
    42 = "Friendly       .           "
    43 = "Lazy           0.1134626736"
    44 = "Disorganized   .           "
    45 = "Unemotional    .           "
    46 = "Imitative      .           "
    47 = "Dominant       0.0344058377"
}

# ---------------------------------------------------------------------------
# 3) New column X: re-emit of the original column A glm() output, with a new
#    "age2" row inserted right after "age" (so everything from row 4 down
#    shifts one row relative to column A).
# ---------------------------------------------------------------------------
$xVals = @{
    1 = "(Intercept)    2.83512517"
    2 = "sex           -1.84693305"
    3 = "age            0.99284123"
    4 = "age2           .         "
    5 = "Outgoing       .         "
    6 = "Helpful        0.17931706"
    7 = "Moody          .         "
    8 = "Organized      .         "
    9 = "Selfconfident  .         "
    10 = "Friendly       0.01925750"
    11 = "Warm           0.09321398"
    12 = "Worrying      -0.28773008"
    13 = "Responsible    .         "
    14 = "Forceful       0.15700508"
    15 = "Lively         .         "
    16 = "Caring         0.18811286"
    17 = "Nervous        .         "
    18 = "Creative       .         "
    19 = "Assertive      0.14280139"
    20 = "Hardworking    .         "
    21 = "Imaginative    .         "
    22 = "Softhearted    .         "
    23 = "Calm           .         "
    24 = "Outspoken      .         "
    25 = "Intelligent    .         "
    26 = "Curious       -0.27351880"
    27 = "Active        -0.12928972"
    28 = "Careless       0.10684492"
    29 = "Broadminded    0.11562209"
    30 = "Sympathetic    .         "
    31 = "Talkative      .         "
    32 = "Sophisticated  0.01474044"
    33 = "Adventurous    .         "
    34 = "Dominant       .         "
    35 = "Thorough       .         "
}

# Apply the existing "Lucida Console" cell format (same as used in columns
# A/E/I/M) to the new S and X ranges before writing values, by copying the
# format from A1 (which already carries that style) via PasteSpecial.
$ws.Range("A1").Copy()
$ws.Range("S1:S47").PasteSpecial(-4122)
$ws.Range("X1:X35").PasteSpecial(-4122)

foreach ($r in $sVals.Keys) {
    $ws.Cells.Item($r, 19).Value = $sVals[$r]
}

foreach ($r in $xVals.Keys) {
    $ws.Cells.Item($r, 24).Value = $xVals[$r]
}

# ---------------------------------------------------------------------------
# 4) Selection moves to X1 (matches the saved sheet view in the target file)
# ---------------------------------------------------------------------------
$ws.Range("X1").Select()
